# Daily attendance processing - 2025-11-12 18:30:39
# Normalizes the "Recorded By" (column G) entries so that any "System"/"system"
# token is moved to the end of the comma-separated list. If no System token is
# present, the first and last entries are swapped instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G$r")
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -lt 2) {
        continue
    }

    $firstIsSystem = ($parts[0].ToLower() -eq "system")

    if ($firstIsSystem) {
        $tmp = $parts[0]
        $parts[0] = $parts[$parts.Count - 1]
        $parts[$parts.Count - 1] = $tmp
        $cell.Value = ($parts -join ", ")
    } else {
        $hasSystemElsewhere = $false
        for ($i = 1; $i -lt $parts.Count; $i++) {
            if ($parts[$i].ToLower() -eq "system") {
                $hasSystemElsewhere = $true
            }
        }

        if (-not $hasSystemElsewhere) {
            $tmp = $parts[0]
            $parts[0] = $parts[$parts.Count - 1]
            $parts[$parts.Count - 1] = $tmp
            $cell.Value = ($parts -join ", ")
        }
    }
}
